$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = 44215
$ws.Range("J2").Value = 16000
$ws.Range("D3").Value2 = 44168
$ws.Range("J3").Value = 7000
$ws.Range("K3").Value = 3000
$ws.Range("M3").Value = 3000
$ws.Range("P3").Value = 30
$ws.Range("D4").Value2 = 44161
$ws.Range("D5").Value2 = 44189
$ws.Range("J5").Value = 16000
$ws.Range("D6").Value2 = 44187
$ws.Range("J6").Value = 12000
$ws.Range("O6").Value = "Provincia de Chacabuco"
$ws.Range("D7").Value2 = 44209
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 7000
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = 2750
$ws.Range("O7").Value = "Provincia de Chacabuco"
$ws.Range("P7").Value = 28
$ws.Range("D9").Value2 = 44245
$ws.Range("J9").Value = 9000
$ws.Range("O9").Value = "Región Metropolitana"
$ws.Range("D10").Value2 = 44245
$ws.Range("I10").Value = "Segunda"
$ws.Range("J10").Value = 5000
$ws.Range("K10").Value = 2500
$ws.Range("L10").Value = 2500
$ws.Range("M10").Value = 2500
$ws.Range("O10").Value = "Región Metropolitana"
$ws.Range("P10").Value = 25
$ws.Range("D12").Value2 = 44159
$ws.Range("D13").Value2 = 44204
$ws.Range("J13").Value = 7000
$ws.Range("D14").Value2 = 44186
$ws.Range("J14").Value = 10000
$ws.Range("D15").Value2 = 44210
$ws.Range("J15").Value = 8800
$ws.Range("K15").Value = 2500
$ws.Range("M15").Value = 2750
$ws.Range("P15").Value = 28
$ws.Range("D16").Value2 = 44188
$ws.Range("J16").Value = 12000
$ws.Range("D17").Value2 = 44166
$ws.Range("D18").Value2 = 44162
$ws.Range("D19").Value2 = 44231
$ws.Range("D20").Value2 = 44232
$ws.Range("J20").Value = 16000
$ws.Range("D21").Value2 = 44229
$ws.Range("J21").Value = 16000
$ws.Range("D22").Value2 = 44214
$ws.Range("D23").Value2 = 44167
$ws.Range("J23").Value = 7000
$ws.Range("D24").Value2 = 44160
$ws.Range("K24").Value = 3000
$ws.Range("M24").Value = 3000
$ws.Range("P24").Value = 30
